$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1:D1").Value = 15
